# Update the Database_Info workbook's "Related Tables" section on Sheet1:
# The Book/Movie "Author/Director" lookup naming was generalized, and a
# new "Makeup_Attributes" row was added at the bottom of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename Author -> Author_Director, Author_Lookup -> AD_Lookup
$ws.Range("H21").Value = "Author_Director"
$ws.Range("H22").Value = "AD_Lookup"

# Rename Color_Lookup -> Feature_Lookup, Size_Lookup -> Format
$ws.Range("H26").Value = "Feature_Lookup"
$ws.Range("H28").Value = "Format"

# Add a new entry at the end of the list
$ws.Range("H31").Value = "Makeup_Attributes"

# Update the view state to match the saved selection/scroll position
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("G26").Select()
